$wb = $excel.ActiveWorkbook

# The localization job moved on from handoff, so the "Status" cells that
# used to read "Ready for handoff" are now "In Translation" everywhere
# that status is reported: the Overview roll-up (one column per language)
# and each language-specific status sheet.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The "Status" columns got narrower now that the longest value in them
# ("Ready for handoff") is gone, so re-fit their width to the new content.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
